# Insert one new weekly price record for "Ajo" (Chino / Primera) at the top
# of the data block (row 139), pushing the existing records (old rows
# 139-196) down by one row to become rows 140-197. This grows the used
# range from A1:R196 to A1:R197.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 139..196 down to 140..197, leaving row 139 blank (but
# inheriting formatting, e.g. the date style on column D) for the new record.
$ws.Rows(139).Insert()

# Populate the newly inserted row 139 with the new record.
$ws.Cells.Item(139, 1).Value = 7
$ws.Cells.Item(139, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(139, 3).Value = "Ñuble"
$ws.Cells.Item(139, 4).Value = 44609
$ws.Cells.Item(139, 5).Value = 16
$ws.Cells.Item(139, 6).Value = 100112003
$ws.Cells.Item(139, 7).Value = "Ajo"
$ws.Cells.Item(139, 8).Value = "Chino"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 60
$ws.Cells.Item(139, 11).Value = 19000
$ws.Cells.Item(139, 12).Value = 20000
$ws.Cells.Item(139, 13).Value = 19500
$ws.Cells.Item(139, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(139, 15).Value = "China"
$ws.Cells.Item(139, 16).Value = 1950
$ws.Cells.Item(139, 17).Value = 10
$ws.Cells.Item(139, 18).Value = "Hortaliza"
